$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Shardul Thakur"

# Insert a new column A (matchNo) before the existing data, shifting
# everything that was there (teamName..result) one column to the right.
$ws.Columns("A:A").Insert()

# Header row
$ws.Cells.Item(1,1).Value = "matchNo"

# Row 2 (existing row, now shifted) gets its matchNo value
$ws.Cells.Item(2,1).Value = "Qualifier"

# New rows of scraped data
$rows = @(
  @("30th","Chennai Super Kings","Shardul Thakur","","1","1","0","0","100.00","Mumbai Indians","Dubai (DSC)","September 19","Super Kings won by 20 runs"),
  @("38th","Chennai Super Kings","Shardul Thakur","","3","2","0","0","150.00","Kolkata Knight Riders","Abu Dhabi","September 26","Super Kings won by 2 wickets"),
  @("12th","Chennai Super Kings","Shardul Thakur","run out (†Samson/Mustafizur Rahman)","1","1","0","0","100.00","Rajasthan Royals","Wankhede","April 19","Super Kings won by 45 runs")
)

$r = 3
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $cell = $ws.Cells.Item($r, $c)
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
      # Force storage as text (matches the source data, which keeps
      # purely-numeric-looking values like "0", "1", "100.00" as text)
      $cell.Value = "'" + $val
    } else {
      $cell.Value = $val
    }
    $c = $c + 1
  }
  $r = $r + 1
}
